$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the "Price" column (D) stays text even for plain-numeric-looking values,
# matching the workbook author convention of storing prices as text.

$ws.Range('D2').Value = '42.999.10'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '2.305.11'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.03'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.55'
$ws.Range('E6').Value = '  -3.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.523'
$ws.Range('E7').Value = '  +3.60%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  +0.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.63'
$ws.Range('E10').Value = '  -1.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0789'
$ws.Range('E11').Value = '  -0.87%  '
$ws.Range('E12').Value = '  -0.94%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '17.89'
$ws.Range('E13').Value = '  -0.34%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.88'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').Value = '2.666.15'
$ws.Range('E15').Value = '  -0.90%  '
$ws.Range('D16').Value = '2.336.51'
$ws.Range('E16').Value = '  +1.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.788'
$ws.Range('D18').Value = '42.926.72'
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.43'
$ws.Range('E19').Value = '  +5.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.18'
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').Value = '0.0₃0909'
$ws.Range('E21').Value = '  +0.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.23'
$ws.Range('E22').Value = '  +0.49%  '
$ws.Range('E23').Value = '  +0.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.16'
$ws.Range('E24').Value = '  -2.88%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.10%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.44'
$ws.Range('E26').Value = '  -1.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.78'
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '168.08'
$ws.Range('E28').Value = '  +0.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.16'
$ws.Range('E29').Value = '  -1.30%  '
$ws.Range('E30').Value = '  -10.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.35'
$ws.Range('E31').Value = '  -3.45%  '
$ws.Range('E32').Value = '  +3.10%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('E34').Value = '  +1.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.23'
$ws.Range('E35').Value = '  +5.32%  '
$ws.Range('E36').Value = '  -0.53%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0692'
$ws.Range('E37').Value = '  -0.21%  '
$ws.Range('E38').Value = '  -1.02%  '
$ws.Range('E39').Value = '  -0.27%  '
$ws.Range('E40').Value = '  +1.20%  '
$ws.Range('E41').Value = '  -2.92%  '
$ws.Range('D42').Value = '1.996.35'
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('E43').Value = '  -0.50%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.08'
$ws.Range('E44').Value = '  -2.19%  '
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.12'
$ws.Range('E45').Value = '  -7.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.47'
$ws.Range('E46').Value = '  -1.53%  '
$ws.Range('E47').Value = '  -2.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.63'
$ws.Range('E48').Value = '  -2.80%  '
$ws.Range('D49').Value = '2.535.44'
$ws.Range('E49').Value = '  +0.54%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.47'
$ws.Range('E50').Value = '  +4.16%  '
$ws.Range('E51').Value = '  -0.19%  '
